# Re-orders the author list in the "Recorded By" column (G) of the
# "Session Analysis Results" sheet.
#
# Sort priority for each comma-separated entry:
#   0 - the exact literal "System"
#   1 - any other casing of "system" (e.g. "system", "SYSTEM")
#   2 - everything else (e-mail addresses), ordered alphabetically
#       (case-insensitive)
# The sort is stable, so equal-priority items keep their relative order
# unless alphabetical ordering (priority 2) dictates otherwise.

function Get-RecorderPriority($item) {
    if ($item.Equals("System")) {
        return 0
    } elseif ($item.ToLower().Equals("system")) {
        return 1
    } else {
        return 2
    }
}

function Compare-Recorders($a, $b) {
    $pa = Get-RecorderPriority $a
    $pb = Get-RecorderPriority $b
    if ($pa -ne $pb) {
        return $pa - $pb
    }
    $la = $a.ToLower()
    $lb = $b.ToLower()
    if ($la.Equals($lb)) {
        return 0
    } elseif ($la -lt $lb) {
        return -1
    } else {
        return 1
    }
}

function Sort-Recorders($items) {
    # Stable insertion sort.
    $n = $items.Count
    for ($i = 1; $i -lt $n; $i++) {
        $key = $items[$i]
        $j = $i - 1
        while (($j -ge 0) -and ((Compare-Recorders $items[$j] $key) -gt 0)) {
            $items[$j + 1] = $items[$j]
            $j = $j - 1
        }
        $items[$j + 1] = $key
    }
    return $items
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($null -eq $val) {
        continue
    }
    if ($val -eq "") {
        continue
    }

    $parts = [string]$val -split ", "
    if ($parts.Count -le 1) {
        continue
    }

    $sorted = Sort-Recorders $parts
    $newVal = $sorted -join ", "

    if (-not $newVal.Equals([string]$val)) {
        $cell.Value2 = $newVal
    }
}

Write-Output "done"
